$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing Effort/Remain values for row 2 (White/Blacklist Implementation)
$ws.Range("D2").Value = "3h"
$ws.Range("E2").Value = "-"

# Fill in Effort/Remain values for row 4 (Rock Paper Scissor against Bot)
$ws.Range("D4").Value = "2h"
$ws.Range("E4").Value = "-"

# Add new row 5 (Coin Flip against User)
$ws.Range("A5").Value = "Coin Flip against User"
$ws.Range("B5").Value = "2h"
$ws.Range("D5").Value = "1h"
$ws.Range("F5").Value = "Tischler, Trinkl"

# Update the selected cell to match the new state
$ws.Range("F6").Select()
